# Updated symbol list on Fri Dec 23 13:48:22 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores its values as plain text (numeric-looking
# strings), so each target cell is first forced to Text format before the
# new value is written -- otherwise Excel would silently re-interpret a
# value like "246.35" as a number. The "Volume(1h)" column (E) holds
# ordinary text and needs no special handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new price text (column D)
$priceUpdates = [ordered]@{
    "D2"  = "246.35"
    "D3"  = "22.11"
    "D4"  = "5.422"
    "D5"  = "0.05878"
    "D6"  = "3.381"
    "D7"  = "6.358"
    "D8"  = "0.8169"
    "D9"  = "1.011"
    "D11" = "0.07466"
    "D12" = "0.03412"
    "D13" = "0.03020"
    "D14" = "4.215"
    "D15" = "0.09384"
    "D16" = "0.001605"
    "D17" = "0.04825"
    "D18" = "0.0005880"
    "D19" = "0.006037"
    "D20" = "0.004125"
    "D21" = "0.0009989"
    "D22" = "0.0001500"
    "D23" = "3.698"
    "D24" = "2.219"
    "D26" = "0.1294"
    "D27" = "0.0001287"
    "D40" = "0.03864"
    "D41" = "0.006620"
    "D42" = "0.1079"
    "D43" = "0.002999"
    "D44" = "0.006235"
    "D45" = "0.00005612"
    "D46" = "0.00000000750"
    "D47" = "0.4989"
    "D48" = "0.1421"
    "D49" = "0.00002100"
    "D50" = "0.01010"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Map of cell address -> new volume/label text (column E)
$textUpdates = [ordered]@{
    "E18" = "17OneONEWorstin24h"
    "E27" = "26UpBotsUBXT"
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

Write-Host "Applied $($priceUpdates.Count) price updates and $($textUpdates.Count) text updates"
